$d = $word.ActiveDocument

# Locate the dependency version string "'com.philips.cdp:uikitLib:1.1.0'"
$full = $d.Content.Text
$uikitIdx = $full.IndexOf("uikitLib:1.1.0")
$verIdx = $full.IndexOf("1.1.0", $uikitIdx)

# Bump the uikitLib dependency version from 1.1.0 to 1.2.0.
$verRange = $d.Range($verIdx, $verIdx + 5)
$verRange.Text = "1.2.0"

# Remove the stale "_GoBack" bookmark that currently sits in its own
# paragraph further down the document.
$d.Bookmarks("_GoBack").Delete()

# Re-split the run right before "cdp:uikitLib" the same way Word leaves
# behind a run boundary at a prior edit point (no visible formatting
# change -- the run properties on both sides stay identical).
$full2 = $d.Content.Text
$cdpIdx = $full2.IndexOf("cdp:uikitLib")
$splitRange = $d.Range($cdpIdx, $cdpIdx)
$d.Bookmarks.Add("TempSplit", $splitRange) | Out-Null
$d.Bookmarks("TempSplit").Delete()

# Re-create "_GoBack" at the new last-edit position, right after the
# "1.2" that was typed (between "1.2" and ".0'").
$bmPos = $verIdx + 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
